$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this pushes the existing rows 8-12 down to 9-13
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with the new record's data
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(8, 3).Value = "Los Lagos"
$ws.Cells.Item(8, 4).Value = 44775
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
$ws.Cells.Item(8, 5).Value = 10
$ws.Cells.Item(8, 6).Value = 100112035
$ws.Cells.Item(8, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 120
$ws.Cells.Item(8, 11).Value = 24000
$ws.Cells.Item(8, 12).Value = 24000
$ws.Cells.Item(8, 13).Value = 24000
$ws.Cells.Item(8, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(8, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(8, 16).Value = 1600
$ws.Cells.Item(8, 17).Value = 15
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Fix Volumen (column J) on what is now row 12 (previously row 11): 80 -> 90
$ws.Cells.Item(12, 10).Value = 90
